$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.48"
$ws.Range("E2").Value = "'0.45%"
$ws.Range("D3").Value = "'32.16"
$ws.Range("E3").Value = "'1.54%"
$ws.Range("D4").Value = "'4.997"
$ws.Range("E4").Value = "'-2.98%"
$ws.Range("D5").Value = "'0.07894"
$ws.Range("E5").Value = "'-2.58%"
$ws.Range("D6").Value = "'2.106"
$ws.Range("E6").Value = "'-17.07%"
$ws.Range("D7").Value = "'7.800"
$ws.Range("E7").Value = "'0.18%"
$ws.Range("D8").Value = "'3.827"
$ws.Range("E8").Value = "'-2.11%"
$ws.Range("D9").Value = "'0.9288"
$ws.Range("E9").Value = "'-0.36%"
$ws.Range("E10").Value = "'-0.50%"
$ws.Range("D11").Value = "'0.07981"
$ws.Range("E11").Value = "'8.54%"
$ws.Range("D12").Value = "'0.08613"
$ws.Range("E12").Value = "'-3.11%"
$ws.Range("D13").Value = "'0.03111"
$ws.Range("E13").Value = "'2.54%"
$ws.Range("E14").Value = "'0.03%"
$ws.Range("D15").Value = "'0.001527"
$ws.Range("E15").Value = "'0.82%"
$ws.Range("D16").Value = "'0.006006"
$ws.Range("E16").Value = "'2.66%"
$ws.Range("E17").Value = "'2,097.54%"
$ws.Range("E18").Value = "'-2.81%"
$ws.Range("E19").Value = "'-0.49%"
$ws.Range("E20").Value = "'0.48%"
$ws.Range("E21").Value = "'-2.37%"
$ws.Range("D22").Value = "'4.280"
$ws.Range("E22").Value = "'2.92%"
$ws.Range("D23").Value = "'0.1792"
$ws.Range("E23").Value = "'6.66%"
$ws.Range("D24").Value = "'0.04599"
$ws.Range("E24").Value = "'-0.72%"
$ws.Range("E25").Value = "'-0.29%"
$ws.Range("D26").Value = "'0.004448"
$ws.Range("D27").Value = "'0.0001252"
$ws.Range("E27").Value = "'4.36%"
$ws.Range("D39").Value = "'0.01717"
$ws.Range("E39").Value = "'-2.52%"
$ws.Range("D40").Value = "'0.04776"
$ws.Range("E40").Value = "'3.59%"
$ws.Range("D41").Value = "'0.007488"
$ws.Range("E41").Value = "'8.14%"
$ws.Range("D42").Value = "'0.1358"
$ws.Range("E42").Value = "'-1.26%"
$ws.Range("D43").Value = "'0.002263"
$ws.Range("E43").Value = "'5.66%"
$ws.Range("D44").Value = "'0.01027"
$ws.Range("E44").Value = "'-0.80%"
$ws.Range("D45").Value = "'0.00005996"
$ws.Range("E45").Value = "'-4.92%"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("D47").Value = "'0.003395"
$ws.Range("E47").Value = "'-59.62%"
$ws.Range("D48").Value = "'0.8204"
$ws.Range("E48").Value = "'9.58%"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("D50").Value = "'0.0002003"
$ws.Range("E50").Value = "'0.05%"
